# Updated symbol list on Thu Dec 15 03:49:50 UTC 2022 with GitHub Actions
#
# Applies the latest coinranking.com price/row refresh to the "cryptos"
# sheet: numeric Price (column D) updates are written back as literal
# text (matching the sheet's existing inlineStr/text convention instead
# of being auto-coerced into floating point numbers), while the
# Coin/Link/Volume columns (B/C/E) for the ProBitToken..BitpandaEcosystemToken
# block are rewritten to reflect ProBitToken's new rank (it moved up to
# position 14, shifting the rows below it down by one).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rangeAddress, $value) {
    # Column D holds price strings like "0.06160" / "0.00006156" that must
    # stay literal text (leading/trailing zeros, no scientific notation).
    # Plain `.Value = "..."` lets Excel's COM layer auto-detect the string
    # as a number, so force the cell to Text format first, then restore
    # the default "Normal" style so no stray formatting is left behind.
    $cell = $ws.Range($rangeAddress)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# ---- Column D (Price) only updates ----
Set-TextValue "D2"  "265.96"
Set-TextValue "D3"  "22.58"
Set-TextValue "D4"  "6.276"
Set-TextValue "D5"  "0.06156"
Set-TextValue "D6"  "3.573"
Set-TextValue "D7"  "6.668"
Set-TextValue "D8"  "1.343"
Set-TextValue "D10" "0.01360"
Set-TextValue "D11" "0.1591"
Set-TextValue "D12" "0.08290"
Set-TextValue "D13" "0.03430"

# ---- Rows 15-26: ProBitToken moves to rank 14, everything below shifts down ----
$ws.Range("B15").Value = "ProBitToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
Set-TextValue "D15" "0.1237"
$ws.Range("E15").Value = "14ProBitTokenPROB"

$ws.Range("B16").Value = "BitMartToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextValue "D16" "0.09258"
$ws.Range("E16").Value = "15BitMartTokenBMX"

$ws.Range("B17").Value = "MCDex"
$ws.Range("C17").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
Set-TextValue "D17" "3.885"
$ws.Range("E17").Value = "16MCDexMCB"

$ws.Range("B18").Value = "BitForexToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextValue "D18" "0.001726"
$ws.Range("E18").Value = "17BitForexTokenBF"

$ws.Range("B19").Value = "CoinExToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
Set-TextValue "D19" "0.04888"
$ws.Range("E19").Value = "18CoinExTokenCET"

$ws.Range("B20").Value = "TigerCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextValue "D20" "0.006217"
$ws.Range("E20").Value = "19TigerCashTCH"

$ws.Range("B21").Value = "HotbitToken"
$ws.Range("C21").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
Set-TextValue "D21" "0.005279"
$ws.Range("E21").Value = "20HotbitTokenHTB"

$ws.Range("B22").Value = "BitKan"
$ws.Range("C22").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
Set-TextValue "D22" "0.001089"
$ws.Range("E22").Value = "21BitKanKAN"

$ws.Range("B23").Value = "NitroEx"
$ws.Range("C23").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
Set-TextValue "D23" "0.0001499"
$ws.Range("E23").Value = "22NitroExNTX"

$ws.Range("B24").Value = "LEO"
$ws.Range("C24").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue "D24" "3.773"
$ws.Range("E24").Value = "23LEOLEO"

$ws.Range("B25").Value = "BTSEToken"
$ws.Range("C25").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
Set-TextValue "D25" "2.276"
$ws.Range("E25").Value = "24BTSETokenBTSE"

$ws.Range("B26").Value = "BitpandaEcosystemToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
Set-TextValue "D26" "0.3341"
$ws.Range("E26").Value = "25BitpandaEcosystemTokenBEST"

# ---- Column D (Price) only updates, further down the sheet ----
Set-TextValue "D40" "0.04616"
Set-TextValue "D41" "0.006959"
Set-TextValue "D42" "0.1137"
Set-TextValue "D43" "0.003399"
Set-TextValue "D44" "0.01079"
Set-TextValue "D45" "0.00006156"
Set-TextValue "D47" "0.6997"
Set-TextValue "D48" "0.1923"
Set-TextValue "D50" "0.01239"
